$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = '2026-02-16T19:15:44+00:00'

# Row 3
$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = '2026-02-16T19:15:45+00:00'

# Row 4
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = '888755853343'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = 'GSU15D03C00NRDP'
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = '4684'
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = '2026-02-15T09:42:36-06:00'
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = 'Racine, WI, 53405, US, United States'
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = '2026-02-16T19:15:42+00:00'
$ws.Range("L4").NumberFormat = "@"
$ws.Range("L4").Value = '888755853343'
$ws.Range("M4").NumberFormat = "@"
$ws.Range("M4").Value = '{"external_order_id": "GSU15D03C00NRDP", "sales_office_id": "4684"}'

# Row 5
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = '888758833785'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = 'CS637989701'
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = '4310'
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = '2026-02-15T12:12:32-05:00'
$ws.Range("J5").NumberFormat = "@"
$ws.Range("J5").Value = 'Lexington, KY, 40514, US, United States'
$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = '2026-02-16T19:15:44+00:00'
$ws.Range("L5").NumberFormat = "@"
$ws.Range("L5").Value = '888758833785'
$ws.Range("M5").NumberFormat = "@"
$ws.Range("M5").Value = '{"external_order_id": "CS637989701", "sales_office_id": "4310"}'

# Row 6
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = '888758749358'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = 'B2BDS10422397'
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = '4461'
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = '2026-02-14T14:26:59-05:00'
$ws.Range("J6").NumberFormat = "@"
$ws.Range("J6").Value = 'Newark, DE, 19702, US, United States'
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = '2026-02-16T19:15:41+00:00'
$ws.Range("L6").NumberFormat = "@"
$ws.Range("L6").Value = '888758749358'
$ws.Range("M6").NumberFormat = "@"
$ws.Range("M6").Value = '{"external_order_id": "B2BDS10422397", "sales_office_id": "4461"}'

# Row 7
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = '888758728566'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = 'PO-211-03682689362551056'
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = '4676'
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = '2026-02-15T10:27:21-06:00'
$ws.Range("J7").NumberFormat = "@"
$ws.Range("J7").Value = 'Pensacola, FL, 32503, US, United States'
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = '2026-02-16T19:15:42+00:00'
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = '888758728566'
$ws.Range("M7").NumberFormat = "@"
$ws.Range("M7").Value = '{"external_order_id": "PO-211-03682689362551056", "sales_office_id": "4676"}'

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = '888755761018'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = 'CS637973380'
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = '4310'
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = '2026-02-15T13:41:01-05:00'
$ws.Range("J8").NumberFormat = "@"
$ws.Range("J8").Value = 'Morristown, TN, 37814, US, United States'
$ws.Range("K8").NumberFormat = "@"
$ws.Range("K8").Value = '2026-02-16T19:14:04+00:00'
$ws.Range("L8").NumberFormat = "@"
$ws.Range("L8").Value = '888755761018'
$ws.Range("M8").NumberFormat = "@"
$ws.Range("M8").Value = '{"external_order_id": "CS637973380", "sales_office_id": "4310"}'

# Row 9
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = '888755101531'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = 'CS637978814'
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = '2026-02-15T15:51:24-05:00'
$ws.Range("J9").NumberFormat = "@"
$ws.Range("J9").Value = 'Sebastian, FL, 32958, US, United States'
$ws.Range("K9").NumberFormat = "@"
$ws.Range("K9").Value = '2026-02-16T19:13:48+00:00'
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = '888755101531'
$ws.Range("M9").NumberFormat = "@"
$ws.Range("M9").Value = '{"external_order_id": "CS637978814", "sales_office_id": "4310"}'

# Row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = '888758647720'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '912003176542533-8545368577'
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = '4362'
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = '2026-02-15T11:48:26-05:00'
$ws.Range("J10").NumberFormat = "@"
$ws.Range("J10").Value = 'Mentor, OH, 44060, US, United States'
$ws.Range("K10").NumberFormat = "@"
$ws.Range("K10").Value = '2026-02-16T19:13:51+00:00'
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = '888758647720'
$ws.Range("M10").NumberFormat = "@"
$ws.Range("M10").Value = '{"external_order_id": "912003176542533-8545368577", "sales_office_id": "4362"}'

# Row 11
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = '888758619065'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '200014488213354'
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = '4260'
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = '2026-02-16T12:19:02-06:00'
$ws.Range("J11").NumberFormat = "@"
$ws.Range("J11").Value = 'Midland, TX, 79707, US, United States'
$ws.Range("K11").NumberFormat = "@"
$ws.Range("K11").Value = '2026-02-16T19:13:59+00:00'
$ws.Range("L11").NumberFormat = "@"
$ws.Range("L11").Value = '888758619065'
$ws.Range("M11").NumberFormat = "@"
$ws.Range("M11").Value = '{"external_order_id": "200014488213354", "sales_office_id": "4260"}'

# Row 12
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = '888758697048'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = 'CS637978814'
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value = '2026-02-15T15:51:24-05:00'
$ws.Range("J12").NumberFormat = "@"
$ws.Range("J12").Value = 'Sebastian, FL, 32958, US, United States'
$ws.Range("K12").NumberFormat = "@"
$ws.Range("K12").Value = '2026-02-16T19:14:05+00:00'
$ws.Range("L12").NumberFormat = "@"
$ws.Range("L12").Value = '888758697048'
$ws.Range("M12").NumberFormat = "@"
$ws.Range("M12").Value = '{"external_order_id": "CS637978814", "sales_office_id": "4310"}'

# Row 13
$ws.Range("K13").NumberFormat = "@"
$ws.Range("K13").Value = '2026-02-16T19:13:49+00:00'

# Row 14
$ws.Range("K14").NumberFormat = "@"
$ws.Range("K14").Value = '2026-02-16T19:13:46+00:00'

# Row 15
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = '888758648050'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = 'CS637973380'
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = '4310'
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = '2026-02-15T13:41:01-05:00'
$ws.Range("J15").NumberFormat = "@"
$ws.Range("J15").Value = 'Morristown, TN, 37814, US, United States'
$ws.Range("K15").NumberFormat = "@"
$ws.Range("K15").Value = '2026-02-16T19:13:48+00:00'
$ws.Range("L15").NumberFormat = "@"
$ws.Range("L15").Value = '888758648050'
$ws.Range("M15").NumberFormat = "@"
$ws.Range("M15").Value = '{"external_order_id": "CS637973380", "sales_office_id": "4310"}'

# Row 16
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = '888755470069'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = 'CS637903827'
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = '4310'
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = '2026-02-16T14:28:04-05:00'
$ws.Range("J16").NumberFormat = "@"
$ws.Range("J16").Value = 'West Branch, MI, 48661, US, United States'
$ws.Range("K16").NumberFormat = "@"
$ws.Range("K16").Value = '2026-02-16T20:14:03+00:00'
$ws.Range("L16").NumberFormat = "@"
$ws.Range("L16").Value = '888755470069'
$ws.Range("M16").NumberFormat = "@"
$ws.Range("M16").Value = '{"external_order_id": "CS637903827", "sales_office_id": "4310"}'

# Row 17
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = '888755594920'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = 'PO-211-03689238948471147'
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = '4676'
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = '2026-02-15T14:25:21-05:00'
$ws.Range("J17").NumberFormat = "@"
$ws.Range("J17").Value = 'Stoughton, MA, 02072, US, United States'
$ws.Range("K17").NumberFormat = "@"
$ws.Range("K17").Value = '2026-02-16T19:13:49+00:00'
$ws.Range("L17").NumberFormat = "@"
$ws.Range("L17").Value = '888755594920'
$ws.Range("M17").NumberFormat = "@"
$ws.Range("M17").Value = '{"external_order_id": "PO-211-03689238948471147", "sales_office_id": "4676"}'

# Row 18
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = '888759745061'
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'fedex'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'FedEx®'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = 'B2BDS10455217'
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = '4461'
$ws.Range("I18").NumberFormat = "@"
$ws.Range("I18").Value = '2026-02-15T10:07:58-05:00'
$ws.Range("J18").NumberFormat = "@"
$ws.Range("J18").Value = 'Windermere, FL, 34786, US, United States'
$ws.Range("K18").NumberFormat = "@"
$ws.Range("K18").Value = '2026-02-16T19:13:46+00:00'
$ws.Range("L18").NumberFormat = "@"
$ws.Range("L18").Value = '888759745061'
$ws.Range("M18").NumberFormat = "@"
$ws.Range("M18").Value = '{"external_order_id": "B2BDS10455217", "sales_office_id": "4461"}'

# Row 19
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = '888757584439'
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'fedex'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'FedEx®'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = 'Delivered'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '115269468-1'
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = '4175'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = 'api'
$ws.Range("H19").Value = ""
$ws.Range("I19").NumberFormat = "@"
$ws.Range("I19").Value = '2026-02-15T12:08:14-05:00'
$ws.Range("J19").NumberFormat = "@"
$ws.Range("J19").Value = 'Virginia Beach, VA, 23453, US, United States'
$ws.Range("K19").NumberFormat = "@"
$ws.Range("K19").Value = '2026-02-16T19:13:47+00:00'
$ws.Range("L19").NumberFormat = "@"
$ws.Range("L19").Value = '888757584439'
$ws.Range("M19").NumberFormat = "@"
$ws.Range("M19").Value = '{"external_order_id": "115269468-1", "sales_office_id": "4175"}'

# Row 20
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = '1036568124332U'
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'dpd-poland'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'DPD Poland'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = 'Delivered'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = 'B2B25338826'
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = '4233'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = 'api'
$ws.Range("H20").Value = ""
$ws.Range("I20").NumberFormat = "@"
$ws.Range("I20").Value = '2026-02-16T08:58:38'
$ws.Range("J20").NumberFormat = "@"
$ws.Range("J20").Value = 'WA2'
$ws.Range("K20").NumberFormat = "@"
$ws.Range("K20").Value = '2026-02-16T13:40:52+00:00'
$ws.Range("L20").NumberFormat = "@"
$ws.Range("L20").Value = '1036568124332U'
$ws.Range("M20").NumberFormat = "@"
$ws.Range("M20").Value = '{"external_order_id": "B2B25338826", "sales_office_id": "4233"}'
